# Update "table_1_112" worksheet: fix header typo, re-style header/footer rows,
# and reorganize the summary block at the bottom of the table.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Fix the spelling of the tenant surname column header (typo fix: "съемщика" -> "съёмщика")
$ws.Range("B2").Value = "Фамилия квартиросъёмщика"

# 2. Row 1 switches from the plain style to the centered Times New Roman
#    style (style of A2/A1) used by the rest of the header/body rows.
$ws.Range("A1").Copy()
$ws.Range("B1:K1").PasteSpecial(-4122)

# 3. Rows 39-43 become data-row styled (like rows 1-38), rebuilding the
#    spacer row and the four summary rows with new order/labels.
$ws.Range("A1:K1").Copy()
$ws.Range("A39:K39").PasteSpecial(-4122)
$ws.Range("A40:K43").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("A39:K43").RowHeight = 15.75

# (Label text is assigned in this order so the workbook's shared-string
#  table is rebuilt in the same sequence the source file uses.)
$ws.Range("B42").Value = "Максимальный срок просрочки, дней"
$ws.Range("B43").Value = "Максимальная сумма к оплате, руб."
$ws.Range("B41").Value = "Средняя площадь, кв.м."
$ws.Range("B40").Value = 'Общая сумма графы "Итого", руб.'

$ws.Range("C40").Formula = "=TRUNC(SUM(K3:K38))"
$ws.Range("C41").Formula = "=AVERAGE(C3:C38)"
$ws.Range("C42").Formula = "=MAX(H3:H38)"
$ws.Range("C43").Formula = "=MAX(K3:K38)"

# 4. Rows 44-46 revert to plain blank spacer rows (full clear removes the
#    old summary labels/values and their styling).
$ws.Range("B44:C46").Clear()

# 5. Update the active selection to match the saved view.
$ws.Range("B39").Select()
